# Mise à jour des résultats du script
# - La cellule C365 ("NA") est vidée.
# - Quatre nouvelles lignes de résultats (366-369) sont ajoutées pour le
#   terme "développement durable" au 2025-12-03.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C365 no longer has a result ("NA" -> blank), D365 (Occurences) stays 1.
$ws.Range("C365").Value = ""

# New rows appended at the bottom of the table.
$newRows = @(
    @{ Row = 366; Date = "2025-12-03"; Terme = "développement durable"; Page = 50; Occ = 1 },
    @{ Row = 367; Date = "2025-12-03"; Terme = "développement durable"; Page = 51; Occ = 1 },
    @{ Row = 368; Date = "2025-12-03"; Terme = "développement durable"; Page = 53; Occ = 1 },
    @{ Row = 369; Date = "2025-12-03"; Terme = "développement durable"; Page = 56; Occ = 1 }
)

foreach ($r in $newRows) {
    # Leading apostrophe forces the date-looking string to stay plain text,
    # matching column A's existing text values (e.g. "2025-12-02").
    $ws.Cells.Item($r.Row, 1).Value = "'" + $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.Terme
    $ws.Cells.Item($r.Row, 3).Value = $r.Page
    $ws.Cells.Item($r.Row, 4).Value = $r.Occ
}
